# Refresh the "cryptos" price/volume table with the latest scrape.
# Generated for commit: "Updated cryptos list ... with GitHub Actions".
#
# Values that look like plain numbers (e.g. "578.19") are written into cells
# that hold them as TEXT (the "Price" column mixes plain decimals with
# dotted-thousands values like "63.980.24", so the whole column is text).
# Force NumberFormat to Text ("@") first so Excel doesn't silently convert
# them to floating point numbers and lose the exact printed digits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.980.24'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '2.752.57'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '578.19'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.03'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('B10').Value = 'TRON'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.384'
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('E12').Value = '  -15.72%  '
$ws.Range('D13').Value = '3.235.53'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.81'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '63.650.10'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('D17').Value = '2.757.98'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.16'
$ws.Range('E18').Value = '  +1.73%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.86'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '358.42'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.82'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.544'
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.53'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.48'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = '0.0₃0905'
$ws.Range('E28').Value = '  -1.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.96'
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.00'
$ws.Range('E30').Value = '  -2.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '171.07'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.21'
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.28'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.93'
$ws.Range('E34').Value = '  +3.41%  '
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.80'
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.983'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.24'
$ws.Range('E39').Value = '  +11.16%  '
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '326.82'
$ws.Range('E41').Value = '  -5.31%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.25'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0592'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.71'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0256'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '135.72'
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.629'
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('E51').Value = '  +0.62%  '
